$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns (I, J), matching the existing header style (s="1")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$srcFormat = $ws.Range("H1")
$srcFormat.Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New numeric data for columns I and J, rows 2-12
$data = @(
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(6, 6),
    @(7, 8),
    @(1, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
